# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same header/style layout used for every
#    per-quarter fund-holdings sheet) and place the copy right before the
#    "总计" (totals) sheet, then rename it "2022-Q1".
# 2. Overwrite the copied sheet's sample row with the real 2022-Q1 holdings
#    (two funds) and add the second data row.
# 3. Update the "总计" sheet: push its existing rows down one and insert the
#    new 2022-Q1 summary row at the top of the data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" sheet from the "2021-Q4" template
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$totals   = $wb.Worksheets.Item("总计")
$template.Copy($totals)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# Worksheet handles returned by Item() are position-based, and inserting the
# new sheet shifted every sheet after it - re-resolve "总计" now that it has
# moved so later writes land on the right physical sheet.
$totals = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# Step 2: write the 2022-Q1 fund holdings onto the new sheet
# ---------------------------------------------------------------------

# Row 2 - 金鹰转型动力灵活配置混合 (overwrite the template's sample row)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "004044"
$newSheet.Range("C2").Value = "金鹰转型动力灵活配置混合"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.72"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.34"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "4.41"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0318"
$newSheet.Range("H2").Value = 8

# Row 3 - 金鹰智慧生活灵活配置混合 (new row; copy A2's style for the index cell)
$newSheet.Range("A2").Copy()
$newSheet.Range("A3").PasteSpecial(-4122)
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "002303"
$newSheet.Range("C3").Value = "金鹰智慧生活灵活配置混合"
$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.11"
$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "89.88"
$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "5.91"
$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0065"
$newSheet.Range("H3").Value = 5

# ---------------------------------------------------------------------
# Step 3: update the "总计" sheet with the 2022-Q1 summary row
# ---------------------------------------------------------------------

# Shift the existing two data rows down by one (row3->row4, row2->row3),
# then write the new 2022-Q1 totals into row2.
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2021-Q3"
$totals.Range("C4").Value = 7
$totals.Range("D4").Value = 0.53

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2021-Q4"
$totals.Range("C3").Value = 1
$totals.Range("D3").Value = 0.02

# A2 already carries the right style (index cell style), reuse it for A4/A3
# too so every index cell in the column matches.
$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)
$totals.Range("A4").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.04
